$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.00581
$ws.Range("H2").Value = 3.01743
$ws.Range("I2").Value = 0.003799625168827527
$ws.Range("J2").Value = 0.003799625168827527
$ws.Range("M2").Value = 29.47402433333333
$ws.Range("N2").Value = 88.422073
$ws.Range("O2").Value = 0.295877356230023
$ws.Range("P2").Value = 0.295877356230023
$ws.Range("Q2").Value = 29.64526841471
$ws.Range("R2").Value = 266.80741573239
$ws.Range("S2").Value = 0.001124223049617743
$ws.Range("T2").Value = 0.001124223049617743

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.00581
$ws.Range("H3").Value = 3.01743
$ws.Range("I3").Value = 0.003799625168827527
$ws.Range("J3").Value = 0.003799625168827527
$ws.Range("O3").Value = 0.1818061388681701
$ws.Range("P3").Value = 0.1818061388681701
$ws.Range("Q3").Value = 18.21596574629
$ws.Range("R3").Value = 163.94369171661
$ws.Range("S3").Value = 0.0006907951810908517
$ws.Range("T3").Value = 0.0006907951810908517

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.00581
$ws.Range("H4").Value = 3.01743
$ws.Range("I4").Value = 0.003799625168827527
$ws.Range("J4").Value = 0.003799625168827527
$ws.Range("M4").Value = 7.238098333333333
$ws.Range("N4").Value = 21.714295
$ws.Range("O4").Value = 0.07266023040422054
$ws.Range("P4").Value = 0.07266023040422054
$ws.Range("Q4").Value = 7.280151684650001
$ws.Range("R4").Value = 65.52136516185
$ws.Range("S4").Value = 0.0002760816402166835
$ws.Range("T4").Value = 0.0002760816402166835

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.00581
$ws.Range("H5").Value = 3.01743
$ws.Range("I5").Value = 0.003799625168827527
$ws.Range("J5").Value = 0.003799625168827527
$ws.Range("M5").Value = 44.79281599999999
$ws.Range("N5").Value = 134.378448
$ws.Range("O5").Value = 0.4496562744975863
$ws.Range("P5").Value = 0.4496562744975864
$ws.Range("Q5").Value = 45.05306226096
$ws.Range("R5").Value = 405.47756034864
$ws.Range("S5").Value = 0.001708525297902248
$ws.Range("T5").Value = 0.001708525297902248

# Row 6
$ws.Range("I6").Value = 0.9594121222074437
$ws.Range("J6").Value = 0.9594121222074438
$ws.Range("M6").Value = 29.47402433333333
$ws.Range("N6").Value = 88.422073
$ws.Range("O6").Value = 0.295877356230023
$ws.Range("P6").Value = 0.295877356230023
$ws.Range("Q6").Value = 7485.483072516531
$ws.Range("R6").Value = 67369.34765264878
$ws.Range("S6").Value = 0.2838683222537742
$ws.Range("T6").Value = 0.2838683222537742

# Row 7
$ws.Range("I7").Value = 0.9594121222074437
$ws.Range("J7").Value = 0.9594121222074438
$ws.Range("O7").Value = 0.1818061388681701
$ws.Range("P7").Value = 0.1818061388681701
$ws.Range("S7").Value = 0.1744270135218523
$ws.Range("T7").Value = 0.1744270135218523

# Row 8
$ws.Range("I8").Value = 0.9594121222074437
$ws.Range("J8").Value = 0.9594121222074438
$ws.Range("M8").Value = 7.238098333333333
$ws.Range("N8").Value = 21.714295
$ws.Range("O8").Value = 0.07266023040422054
$ws.Range("P8").Value = 0.07266023040422054
$ws.Range("Q8").Value = 1838.25126622094
$ws.Range("R8").Value = 16544.26139598845
$ws.Range("S8").Value = 0.06971110585219505
$ws.Range("T8").Value = 0.06971110585219506

# Row 9
$ws.Range("I9").Value = 0.9594121222074437
$ws.Range("J9").Value = 0.9594121222074438
$ws.Range("M9").Value = 44.79281599999999
$ws.Range("N9").Value = 134.378448
$ws.Range("O9").Value = 0.4496562744975863
$ws.Range("P9").Value = 0.4496562744975864
$ws.Range("Q9").Value = 11375.97845975679
$ws.Range("R9").Value = 102383.8061378112
$ws.Range("S9").Value = 0.4314056805796221
$ws.Range("T9").Value = 0.4314056805796223

# Row 10
$ws.Range("G10").Value = 9.336668333333334
$ws.Range("H10").Value = 28.010005
$ws.Range("I10").Value = 0.03527091597053946
$ws.Range("J10").Value = 0.03527091597053946
$ws.Range("M10").Value = 29.47402433333333
$ws.Range("N10").Value = 88.422073
$ws.Range("O10").Value = 0.295877356230023
$ws.Range("P10").Value = 0.295877356230023
$ws.Range("Q10").Value = 275.1891896489295
$ws.Range("R10").Value = 2476.702706840365
$ws.Range("S10").Value = 0.01043586536917451
$ws.Range("T10").Value = 0.01043586536917451

# Row 11
$ws.Range("G11").Value = 9.336668333333334
$ws.Range("H11").Value = 28.010005
$ws.Range("I11").Value = 0.03527091597053946
$ws.Range("J11").Value = 0.03527091597053946
$ws.Range("O11").Value = 0.1818061388681701
$ws.Range("P11").Value = 0.1818061388681701
$ws.Range("Q11").Value = 169.0939944367928
$ws.Range("R11").Value = 1521.845949931135
$ws.Range("S11").Value = 0.006412469046947457
$ws.Range("T11").Value = 0.006412469046947457

# Row 12
$ws.Range("G12").Value = 9.336668333333334
$ws.Range("H12").Value = 28.010005
$ws.Range("I12").Value = 0.03527091597053946
$ws.Range("J12").Value = 0.03527091597053946
$ws.Range("M12").Value = 7.238098333333333
$ws.Range("N12").Value = 21.714295
$ws.Range("O12").Value = 0.07266023040422054
$ws.Range("P12").Value = 0.07266023040422054
$ws.Range("Q12").Value = 67.57972350238612
$ws.Range("R12").Value = 608.217511521475
$ws.Range("S12").Value = 0.002562792880987299
$ws.Range("T12").Value = 0.002562792880987299

# Row 13
$ws.Range("G13").Value = 9.336668333333334
$ws.Range("H13").Value = 28.010005
$ws.Range("I13").Value = 0.03527091597053946
$ws.Range("J13").Value = 0.03527091597053946
$ws.Range("M13").Value = 44.79281599999999
$ws.Range("N13").Value = 134.378448
$ws.Range("O13").Value = 0.4496562744975863
$ws.Range("P13").Value = 0.4496562744975864
$ws.Range("Q13").Value = 418.2156667080266
$ws.Range("R13").Value = 3763.94100037224
$ws.Range("S13").Value = 0.01585978867343019
$ws.Range("T13").Value = 0.0158597886734302

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.4016586666666667
$ws.Range("H14").Value = 1.204976
$ws.Range("I14").Value = 0.001517336653189343
$ws.Range("J14").Value = 0.001517336653189343
$ws.Range("M14").Value = 29.47402433333333
$ws.Range("N14").Value = 88.422073
$ws.Range("O14").Value = 0.295877356230023
$ws.Range("P14").Value = 0.295877356230023
$ws.Range("Q14").Value = 11.83849731502755
$ws.Range("R14").Value = 106.546475835248
$ws.Range("S14").Value = 0.000448945557456574
$ws.Range("T14").Value = 0.000448945557456574

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.4016586666666667
$ws.Range("H15").Value = 1.204976
$ws.Range("I15").Value = 0.001517336653189343
$ws.Range("J15").Value = 0.001517336653189343
$ws.Range("O15").Value = 0.1818061388681701
$ws.Range("P15").Value = 0.1818061388681701
$ws.Range("Q15").Value = 7.274336617950222
$ws.Range("R15").Value = 65.46902956155201
$ws.Range("S15").Value = 0.0002758611182795061
$ws.Range("T15").Value = 0.0002758611182795061

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.4016586666666667
$ws.Range("H16").Value = 1.204976
$ws.Range("I16").Value = 0.001517336653189343
$ws.Range("J16").Value = 0.001517336653189343
$ws.Range("M16").Value = 7.238098333333333
$ws.Range("N16").Value = 21.714295
$ws.Range("O16").Value = 0.07266023040422054
$ws.Range("P16").Value = 0.07266023040422054
$ws.Range("Q16").Value = 2.907244925768889
$ws.Range("R16").Value = 26.16520433192
$ws.Range("S16").Value = 0.0001102500308215065
$ws.Range("T16").Value = 0.0001102500308215065

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.4016586666666667
$ws.Range("H17").Value = 1.204976
$ws.Range("I17").Value = 0.001517336653189343
$ws.Range("J17").Value = 0.001517336653189343
$ws.Range("M17").Value = 44.79281599999999
$ws.Range("N17").Value = 134.378448
$ws.Range("O17").Value = 0.4496562744975863
$ws.Range("P17").Value = 0.4496562744975864
$ws.Range("Q17").Value = 17.99142275080533
$ws.Range("R17").Value = 161.922804757248
$ws.Range("S17").Value = 0.000682279946631756
$ws.Range("T17").Value = 0.000682279946631756
